$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 348; this shifts the existing rows 348:373
# down to 349:374 (values, styles and all), matching the target diff where
# a brand-new weekly price record is inserted ahead of the previously-last
# record and everything below cascades down by one row.
$ws.Rows.Item(348).Insert()

# Populate the newly inserted row 348 with the new weekly record.
$ws.Range("A348").Value = 10
$ws.Range("B348").Value = "Vega Modelo de Temuco"
$ws.Range("C348").Value = "La Araucanía"
$ws.Range("D348").Value = 44746
$ws.Range("E348").Value = 9
$ws.Range("F348").Value = 100112037
$ws.Range("G348").Value = "Cebollín"
$ws.Range("H348").Value = "Sin especificar"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 110
$ws.Range("K348").Value = 9000
$ws.Range("L348").Value = 9000
$ws.Range("M348").Value = 9000
$ws.Range("N348").Value = "$/docena de paquetes"
$ws.Range("O348").Value = "Provincia de Cautín"
$ws.Range("P348").Value = 750
$ws.Range("Q348").Value = 12
$ws.Range("R348").Value = "Hortaliza"

# Match the D column's date number format used by the rest of the sheet.
$ws.Range("D348").NumberFormat = $ws.Range("D349").NumberFormat
